$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 19 de Octubre de 2020 a las 13:07"

# Row 16: Iran -> Iran
$ws.Cells.Item(16, 2).Value = 534631
$ws.Cells.Item(16, 3).Value = 4251
$ws.Cells.Item(16, 4).Value = 431360
$ws.Cells.Item(16, 5).Value = 72559
$ws.Cells.Item(16, 6).Value = 0
$ws.Cells.Item(16, 7).Value = 337
$ws.Cells.Item(16, 8).Value = 30712

# Row 20: Banglades -> Banglades
$ws.Cells.Item(20, 2).Value = 390206
$ws.Cells.Item(20, 3).Value = 1637
$ws.Cells.Item(20, 4).Value = 305599
$ws.Cells.Item(20, 5).Value = 78926
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(20, 7).Value = 21
$ws.Cells.Item(20, 8).Value = 5681

# Row 33: Rumania -> Rumania
$ws.Cells.Item(33, 2).Value = 182854
$ws.Cells.Item(33, 3).Value = 2466
$ws.Cells.Item(33, 4).Value = 132082
$ws.Cells.Item(33, 5).Value = 44841
$ws.Cells.Item(33, 6).Value = 0
$ws.Cells.Item(33, 7).Value = 59
$ws.Cells.Item(33, 8).Value = 5931

# Row 38: Nepal -> Nepal
$ws.Cells.Item(38, 2).Value = 136036
$ws.Cells.Item(38, 3).Value = 3790
$ws.Cells.Item(38, 4).Value = 94501
$ws.Cells.Item(38, 5).Value = 40778
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(38, 7).Value = 18
$ws.Cells.Item(38, 8).Value = 757

# Row 54: Bielorrusia -> Bielorrusia
$ws.Cells.Item(54, 2).Value = 88290
$ws.Cells.Item(54, 3).Value = 592
$ws.Cells.Item(54, 4).Value = 80130
$ws.Cells.Item(54, 5).Value = 7227
$ws.Cells.Item(54, 6).Value = 0
$ws.Cells.Item(54, 7).Value = 4
$ws.Cells.Item(54, 8).Value = 933

# Row 57: Barein -> Suiza
$ws.Cells.Item(57, 1).Value = "Suiza"
$ws.Cells.Item(57, 2).Value = 83159
$ws.Cells.Item(57, 3).Value = 8737
$ws.Cells.Item(57, 4).Value = 50600
$ws.Cells.Item(57, 5).Value = 30427
$ws.Cells.Item(57, 6).Value = 0
$ws.Cells.Item(57, 7).Value = 9
$ws.Cells.Item(57, 8).Value = 2132

# Row 58: Suiza -> Barein
$ws.Cells.Item(58, 1).Value = "Barein"
$ws.Cells.Item(58, 2).Value = 77902
$ws.Cells.Item(58, 3).Value = 0
$ws.Cells.Item(58, 4).Value = 74320
$ws.Cells.Item(58, 5).Value = 3282
$ws.Cells.Item(58, 6).Value = 0
$ws.Cells.Item(58, 7).Value = 0
$ws.Cells.Item(58, 8).Value = 300

# Row 72: Ghana -> Estado de Palestina
$ws.Cells.Item(72, 1).Value = "Estado de Palestina"
$ws.Cells.Item(72, 2).Value = 47616
$ws.Cells.Item(72, 3).Value = 481
$ws.Cells.Item(72, 4).Value = 40861
$ws.Cells.Item(72, 5).Value = 6342
$ws.Cells.Item(72, 6).Value = 0
$ws.Cells.Item(72, 7).Value = 5
$ws.Cells.Item(72, 8).Value = 413

# Row 73: Estado de Palestina -> Ghana
$ws.Cells.Item(73, 1).Value = "Ghana"
$ws.Cells.Item(73, 2).Value = 47372
$ws.Cells.Item(73, 3).Value = 62
$ws.Cells.Item(73, 4).Value = 46664
$ws.Cells.Item(73, 5).Value = 398
$ws.Cells.Item(73, 6).Value = 0
$ws.Cells.Item(73, 7).Value = 0
$ws.Cells.Item(73, 8).Value = 310

# Row 92: Malasia -> Malasia
$ws.Cells.Item(92, 2).Value = 21363
$ws.Cells.Item(92, 3).Value = 865
$ws.Cells.Item(92, 4).Value = 13717
$ws.Cells.Item(92, 5).Value = 7456
$ws.Cells.Item(92, 6).Value = 0
$ws.Cells.Item(92, 7).Value = 3
$ws.Cells.Item(92, 8).Value = 190

# Row 100: Senegal -> Senegal
$ws.Cells.Item(100, 2).Value = 15432
$ws.Cells.Item(100, 3).Value = 14
$ws.Cells.Item(100, 4).Value = 13865
$ws.Cells.Item(100, 5).Value = 1248
$ws.Cells.Item(100, 6).Value = 0
$ws.Cells.Item(100, 7).Value = 2
$ws.Cells.Item(100, 8).Value = 319

# Row 102: Finlandia -> Eslovenia
$ws.Cells.Item(102, 1).Value = "Eslovenia"
$ws.Cells.Item(102, 2).Value = 13679
$ws.Cells.Item(102, 3).Value = 537
$ws.Cells.Item(102, 4).Value = 6385
$ws.Cells.Item(102, 5).Value = 7104
$ws.Cells.Item(102, 6).Value = 0
$ws.Cells.Item(102, 7).Value = 2
$ws.Cells.Item(102, 8).Value = 190

# Row 103: Eslovenia -> Finlandia
$ws.Cells.Item(103, 1).Value = "Finlandia"
$ws.Cells.Item(103, 2).Value = 13555
$ws.Cells.Item(103, 3).Value = 131
$ws.Cells.Item(103, 4).Value = 9100
$ws.Cells.Item(103, 5).Value = 4104
$ws.Cells.Item(103, 6).Value = 0
$ws.Cells.Item(103, 7).Value = 0
$ws.Cells.Item(103, 8).Value = 351

# Row 126: Sri Lanka -> Sri Lanka
$ws.Cells.Item(126, 2).Value = 5585
$ws.Cells.Item(126, 3).Value = 47
$ws.Cells.Item(126, 4).Value = 3440
$ws.Cells.Item(126, 5).Value = 2132
$ws.Cells.Item(126, 6).Value = 0
$ws.Cells.Item(126, 7).Value = 0
$ws.Cells.Item(126, 8).Value = 13

# Row 139: Malta -> Malta
$ws.Cells.Item(139, 2).Value = 4737
$ws.Cells.Item(139, 3).Value = 109
$ws.Cells.Item(139, 4).Value = 3242
$ws.Cells.Item(139, 5).Value = 1450
$ws.Cells.Item(139, 6).Value = 0
$ws.Cells.Item(139, 7).Value = 0
$ws.Cells.Item(139, 8).Value = 45

# Row 176: Gibraltar -> Gibraltar
$ws.Cells.Item(176, 2).Value = 577
$ws.Cells.Item(176, 3).Value = 6
$ws.Cells.Item(176, 4).Value = 462
$ws.Cells.Item(176, 5).Value = 115
$ws.Cells.Item(176, 6).Value = 0
$ws.Cells.Item(176, 7).Value = 0
$ws.Cells.Item(176, 8).Value = 0

# Row 191: Liechtenstein -> Liechtenstein
$ws.Cells.Item(191, 2).Value = 224
$ws.Cells.Item(191, 3).Value = 0
$ws.Cells.Item(191, 4).Value = 142
$ws.Cells.Item(191, 5).Value = 81
$ws.Cells.Item(191, 6).Value = 0
$ws.Cells.Item(191, 7).Value = 0
$ws.Cells.Item(191, 8).Value = 1
